$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2255892255892256
$ws.Range("C2").Value = 0.4915824915824916
$ws.Range("J2").Value = 0.01683501683501683
$ws.Range("P2").Value = 0.1582491582491583
$ws.Range("S2").Value = 0.1077441077441077
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.0245398773006135
$ws.Range("J3").Value = 0.049079754601227
$ws.Range("P3").Value = 0.7852760736196319
$ws.Range("S3").Value = 0.1349693251533742
$ws.Range("J4").Value = 0.09523809523809523
$ws.Range("P4").Value = 0.5476190476190477
$ws.Range("S4").Value = 0.3571428571428572
$ws.Range("B6").Value = 0.06310679611650485
$ws.Range("D6").Value = 0.009708737864077669
$ws.Range("F6").Value = 0.03398058252427184
$ws.Range("J6").Value = 0.3155339805825243
$ws.Range("O6").Value = 0.02427184466019417
$ws.Range("Q6").Value = 0.1553398058252427
$ws.Range("R6").Value = 0.07766990291262135
$ws.Range("S6").Value = 0.3203883495145631
$ws.Range("B7").Value = 0.125
$ws.Range("D7").Value = 0.02840909090909091
$ws.Range("F7").Value = 0.05113636363636364
$ws.Range("J7").Value = 0.1306818181818182
$ws.Range("O7").Value = 0.005681818181818182
$ws.Range("Q7").Value = 0.1193181818181818
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.4772727272727273
$ws.Range("B8").Value = 0.102974828375286
$ws.Range("D8").Value = 0.01601830663615561
$ws.Range("F8").Value = 0.06407322654462243
$ws.Range("J8").Value = 0.1121281464530892
$ws.Range("O8").Value = 0.0137299771167048
$ws.Range("Q8").Value = 0.1601830663615561
$ws.Range("R8").Value = 0.08695652173913043
$ws.Range("S8").Value = 0.4439359267734554
$ws.Range("B9").Value = 0.07199999999999999
$ws.Range("D9").Value = 0.008
$ws.Range("F9").Value = 0.096
$ws.Range("J9").Value = 0.16
$ws.Range("O9").Value = 0.032
$ws.Range("Q9").Value = 0.136
$ws.Range("R9").Value = 0.096
$ws.Range("S9").Value = 0.4
$ws.Range("B10").Value = 0.1170858629661752
$ws.Range("D10").Value = 0.02341717259323504
$ws.Range("E10").Value = 0.002601908065915004
$ws.Range("F10").Value = 0.07545533391153512
$ws.Range("J10").Value = 0.1274934952298352
$ws.Range("O10").Value = 0.02428447528187338
$ws.Range("Q10").Value = 0.2116218560277537
$ws.Range("R10").Value = 0.08326105810928014
$ws.Range("S10").Value = 0.3347788378143972
$ws.Range("F11").Value = 0.00819672131147541
$ws.Range("G11").Value = 0.1475409836065574
$ws.Range("J11").Value = 0.0778688524590164
$ws.Range("K11").Value = 0.1844262295081967
$ws.Range("L11").Value = 0.569672131147541
$ws.Range("S11").Value = 0.01229508196721311
$ws.Range("G12").Value = 0.7517241379310344
$ws.Range("J12").Value = 0.1793103448275862
$ws.Range("K12").Value = 0.006896551724137931
$ws.Range("L12").Value = 0.04137931034482759
$ws.Range("S12").Value = 0.02068965517241379
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.07272727272727272
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02105263157894737
$ws.Range("H15").Value = 0.1526315789473684
$ws.Range("I15").Value = 0.06315789473684211
$ws.Range("J15").Value = 0.3368421052631579
$ws.Range("K15").Value = 0.06315789473684211
$ws.Range("M15").Value = 0.01052631578947368
$ws.Range("O15").Value = 0.04210526315789474
$ws.Range("S15").Value = 0.3105263157894737
$ws.Range("F16").Value = 0.02631578947368421
$ws.Range("H16").Value = 0.1526315789473684
$ws.Range("I16").Value = 0.04736842105263158
$ws.Range("J16").Value = 0.4210526315789473
$ws.Range("K16").Value = 0.1578947368421053
$ws.Range("M16").Value = 0.02105263157894737
$ws.Range("O16").Value = 0.04736842105263158
$ws.Range("S16").Value = 0.1263157894736842
$ws.Range("F17").Value = 0.01275510204081633
$ws.Range("H17").Value = 0.1760204081632653
$ws.Range("I17").Value = 0.06887755102040816
$ws.Range("J17").Value = 0.4158163265306122
$ws.Range("K17").Value = 0.08928571428571429
$ws.Range("M17").Value = 0.03061224489795918
$ws.Range("N17").Value = 0.002551020408163265
$ws.Range("O17").Value = 0.07397959183673469
$ws.Range("S17").Value = 0.1301020408163265
$ws.Range("F18").Value = 0.005617977528089887
$ws.Range("H18").Value = 0.1573033707865168
$ws.Range("I18").Value = 0.05617977528089887
$ws.Range("J18").Value = 0.4831460674157304
$ws.Range("K18").Value = 0.09550561797752809
$ws.Range("M18").Value = 0.02247191011235955
$ws.Range("N18").Value = 0.005617977528089887
$ws.Range("O18").Value = 0.03932584269662921
$ws.Range("S18").Value = 0.1348314606741573
$ws.Range("F19").Value = 0.01612903225806452
$ws.Range("H19").Value = 0.2544802867383513
$ws.Range("I19").Value = 0.06093189964157706
$ws.Range("J19").Value = 0.3593189964157706
$ws.Range("K19").Value = 0.09677419354838709
$ws.Range("M19").Value = 0.02956989247311828
$ws.Range("N19").Value = 0.0008960573476702509
$ws.Range("O19").Value = 0.07168458781362007
$ws.Range("S19").Value = 0.1102150537634409
